# Auto-generated edit script applying the scheduled-runner price update
# described in the commit diff. Each block updates the H-N price/profit
# columns for a specific (sheet, row) pair.
$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1199.4445
$ws.Range("I18").Value = 1199.4445
$ws.Range("K18").Value = 1199.4445
$ws.Range("M18").Value = -915.4445000000001

# ALC row 70
$ws.Range("H70").Value = 2230.8
$ws.Range("I70").Value = 1628.5714
$ws.Range("K70").Value = 4885.7142
$ws.Range("M70").Value = -4615.7142

# ALC row 73
$ws.Range("H73").Value = 2230.8
$ws.Range("I73").Value = 1628.5714
$ws.Range("K73").Value = 4885.7142
$ws.Range("M73").Value = -3949.7142

# ALC row 86
$ws.Range("H86").Value = 22121.309
$ws.Range("J86").Value = 43296.332
$ws.Range("L86").Value = 43296.332
$ws.Range("N86").Value = -45542.332

# ALC row 89
$ws.Range("H89").Value = 22121.309
$ws.Range("J89").Value = 43296.332
$ws.Range("L89").Value = 216481.66
$ws.Range("N89").Value = -227713.66

# ALC row 98
$ws.Range("H98").Value = 2880.9443
$ws.Range("I98").Value = 3047.3125
$ws.Range("J98").Value = 1550
$ws.Range("K98").Value = 3047.3125
$ws.Range("L98").Value = 1550
$ws.Range("M98").Value = -1549.3125
$ws.Range("N98").Value = -4546

# ALC row 111
$ws.Range("H111").Value = 2221.4285
$ws.Range("I111").Value = 1975
$ws.Range("J111").Value = 2320
$ws.Range("K111").Value = 5925
$ws.Range("L111").Value = 6960
$ws.Range("M111").Value = -2858
$ws.Range("N111").Value = -13094

# ALC row 122
$ws.Range("H122").Value = 2880.9443
$ws.Range("I122").Value = 3047.3125
$ws.Range("J122").Value = 1550
$ws.Range("K122").Value = 9141.9375
$ws.Range("L122").Value = 4650
$ws.Range("M122").Value = -6691.9375
$ws.Range("N122").Value = -9550

# ALC row 132
$ws.Range("H132").Value = 12990.808
$ws.Range("I132").Value = 11918.223
$ws.Range("J132").Value = 15404.125
$ws.Range("K132").Value = 35754.669
$ws.Range("L132").Value = 46212.375
$ws.Range("M132").Value = -33224.669
$ws.Range("N132").Value = -51272.375

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2987
$ws.Range("I88").Value = 2250
$ws.Range("J88").Value = 3601.1667
$ws.Range("K88").Value = 2250
$ws.Range("L88").Value = 3601.1667
$ws.Range("M88").Value = -1844
$ws.Range("N88").Value = -4413.1667

# ARM row 91
$ws.Range("H91").Value = 2987
$ws.Range("I91").Value = 2250
$ws.Range("J91").Value = 3601.1667
$ws.Range("K91").Value = 2250
$ws.Range("L91").Value = 3601.1667
$ws.Range("M91").Value = -846
$ws.Range("N91").Value = -6409.1667

# BSM row 11
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 396.66666
$ws.Range("I11").Value = 227.5
$ws.Range("J11").Value = 735
$ws.Range("K11").Value = 227.5
$ws.Range("L11").Value = 735
$ws.Range("M11").Value = -87.5
$ws.Range("N11").Value = -1015

# BSM row 99
$ws.Range("H99").Value = 14989.6
$ws.Range("I99").Value = 1355.75
$ws.Range("J99").Value = 30571.143
$ws.Range("K99").Value = 1355.75
$ws.Range("L99").Value = 30571.143
$ws.Range("M99").Value = 142.25
$ws.Range("N99").Value = -33567.143

# BSM row 134
$ws.Range("H134").Value = 5862.44
$ws.Range("I134").Value = 1583.6888
$ws.Range("K134").Value = 4751.0664
$ws.Range("M134").Value = -2216.0664

# CRP row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 24950
$ws.Range("J95").Value = 24950
$ws.Range("L95").Value = 24950
$ws.Range("N95").Value = -30442

# CRP row 125
$ws.Range("H125").Value = 83998.625
$ws.Range("J125").Value = 83998.625
$ws.Range("L125").Value = 83998.625
$ws.Range("N125").Value = -88918.625

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 138987.5
$ws.Range("J37").Value = 138987.5
$ws.Range("L37").Value = 416962.5
$ws.Range("N37").Value = -417186.5

# CUL row 107
$ws.Range("H107").Value = 7815499.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 7815499.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 23446498.5
$ws.Range("N107").Value = -23450338.5
$ws.Range("M107").ClearContents()

# CUL row 137
$ws.Range("H137").Value = 1429.8572
$ws.Range("I137").Value = 1021.5
$ws.Range("J137").Value = 1974.3334
$ws.Range("K137").Value = 3064.5
$ws.Range("L137").Value = 5923.0002
$ws.Range("M137").Value = 2035.5
$ws.Range("N137").Value = -16123.0002

# GSM row 20
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 12631.667

# GSM row 101
$ws.Range("H101").Value = 50657
$ws.Range("J101").Value = 50657
$ws.Range("L101").Value = 50657
$ws.Range("N101").Value = -57147

# GSM row 102
$ws.Range("H102").Value = 5964
$ws.Range("I102").Value = 6360.4
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 6360.4
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -4738.4
$ws.Range("N102").Value = -5244

# GSM row 108
$ws.Range("H108").Value = 49975
$ws.Range("J108").Value = 49975
$ws.Range("L108").Value = 49975
$ws.Range("N108").Value = -57655

# GSM row 111
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -46134

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1675.2
$ws.Range("J46").Value = 2015.125
$ws.Range("L46").Value = 2015.125
$ws.Range("N46").Value = -2391.125

# LTW row 55
$ws.Range("H55").Value = 1951.6976
$ws.Range("I55").Value = 923.125
$ws.Range("J55").Value = 3250.9473
$ws.Range("K55").Value = 923.125
$ws.Range("L55").Value = 3250.9473
$ws.Range("M55").Value = -750.125
$ws.Range("N55").Value = -3596.9473

# LTW row 100
$ws.Range("H100").Value = 6295.222
$ws.Range("I100").Value = 3506
$ws.Range("K100").Value = 3506
$ws.Range("M100").Value = -2965

# LTW row 103
$ws.Range("H103").Value = 27249.75
$ws.Range("J103").Value = 27249.75
$ws.Range("L103").Value = 27249.75
$ws.Range("N103").Value = -29593.75

# LTW row 132
$ws.Range("H132").Value = 2119032
$ws.Range("I132").Value = 3263.818
$ws.Range("J132").Value = 5028213
$ws.Range("K132").Value = 9791.454000000002
$ws.Range("L132").Value = 15084639
$ws.Range("M132").Value = -7261.454000000002
$ws.Range("N132").Value = -15089699

# LTW row 136
$ws.Range("H136").Value = 18758.5
$ws.Range("I136").Value = 25109.889
$ws.Range("K136").Value = 75329.667
$ws.Range("M136").Value = -72779.667

# WVR row 20
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 14997.5
$ws.Range("I20").Value = 19995
$ws.Range("K20").Value = 19995
$ws.Range("M20").Value = -19755

# WVR row 31
$ws.Range("H31").Value = 17000
$ws.Range("J31").Value = 20000
$ws.Range("L31").Value = 20000
$ws.Range("N31").Value = -20696

# WVR row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

# WVR row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

# WVR row 132
$ws.Range("H132").Value = 8281
$ws.Range("I132").Value = 4676.968
$ws.Range("J132").Value = 17591.416
$ws.Range("K132").Value = 14030.904
$ws.Range("L132").Value = 52774.24800000001
$ws.Range("M132").Value = -11500.904
$ws.Range("N132").Value = -57834.24800000001

# WVR row 136
$ws.Range("H136").Value = 15740.823
$ws.Range("I136").Value = 3157
$ws.Range("J136").Value = 20984.084
$ws.Range("K136").Value = 9471
$ws.Range("L136").Value = 62952.25199999999
$ws.Range("M136").Value = -6921
$ws.Range("N136").Value = -68052.25199999999

